$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add Sheet2, Sheet3, Sheet4 right after Sheet1, in order. Copying Sheet1
# (instead of Worksheets.Add()) keeps the same worksheet markup/namespaces
# that the workbook already uses.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Sheet3"

$ws3.Copy($null, $ws3)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "Sheet4"

# Sheet1 now holds the plain number 1 (the text that used to live there
# moved over to Sheet2).
$ws1.Range("A1").Value = 1

# Sheet2 gets the text value (previously "This is neat!", now "undefined").
$ws2.Range("A1").Value = "undefined"

# Sheet3 / Sheet4 each hold a simple number matching their sheet index.
$ws3.Range("A1").Value = 3
$ws4.Range("A1").Value = 4

# Give Sheet1's tab an accent color (theme Accent2, Darker 25%).
$ws1.Tab.Color = 1137349

# Sheet4 keeps cell E18 selected.
$null = $ws4.Range("E18").Select()

# Keep Sheet1 as the active/selected tab with A1 selected.
$ws1.Activate()
$null = $ws1.Range("A1").Select()

Write-Output "done"
